$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Relocate the existing "COSTO OPERATIVO / TOTAL GASTO / UTILIDAD FINAL" headers
#    three columns to the right (AS8:AU8), preserving their currency-header format.
$ws.Range("AP8").Copy()
$ws.Range("AS8:AU8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AS8").Value = "COSTO OPERATIVO"
$ws.Range("AT8").Value = "TOTAL GASTO"
$ws.Range("AU8").Value = "UTILIDAD FINAL"

# 2) Overwrite the original cells with the new header text (format is already correct)
$ws.Range("AP8").Value = "SUB TOTAL"
$ws.Range("AQ8").Value = "TOTAL"
$ws.Range("AR8").Value = "UTILIDAD BRUTA"

Write-Output "done"
